# WorkflowInput.xlsx update:
# Introduce a dedicated "tpDictionaryLoq" sheet (a duplicate of "tpDictionary")
# and move the "lloq" mapping row out of "tpDictionary" into it.

$wb = $excel.ActiveWorkbook

# Duplicate "tpDictionary" -> new sheet lands immediately after it.
$tpDictionary = $wb.Worksheets.Item("tpDictionary")
$tpDictionary.Copy($null, $tpDictionary)

$tpDictionaryLoq = $wb.Worksheets.Item($tpDictionary.Index + 1)
$tpDictionaryLoq.Name = "tpDictionaryLoq"

# The copied sheet keeps the wrapped-text styling from the original; the
# published version drops word-wrap on the duplicated sheet.
$tpDictionaryLoq.UsedRange.WrapText = $false

# "tpDictionary" keeps only the generic identifier rows; the lloq/LOQ row
# (row 12: lloq / timeprofile / LOQ / ... / units are defined in
# corresponding output definitions) now lives solely on tpDictionaryLoq.
$tpDictionary.Rows.Item(12).Delete()

# Leave the new sheet active/selected, matching the authored workbook view.
[void]$tpDictionaryLoq.Activate()
[void]$tpDictionaryLoq.Range("C15").Select()
[void]$tpDictionary.Range("C17").Select()
[void]$tpDictionaryLoq.Activate()
